# ---------------------------------------------------------------------------
# Applies the "started properties of the data" edit to report.docx
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

function FindReplace($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "FindReplace failed for: $old"
    }
    return $ok
}

# Inserts a new paragraph (Body style / author's own runs inherit formatting
# from the paragraph that is split) right after the given anchor text, and
# types $text into that freshly created paragraph.
function InsertParagraphAfterAnchor($anchor, $text) {
    $r = $d.Content
    $ok = $r.Find.Execute($anchor)
    if (-not $ok) {
        throw "InsertParagraphAfterAnchor: anchor not found: $anchor"
    }
    $endPos = $r.End
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newPos = $endPos + 1
    $ins = $d.Range($newPos, $newPos)
    $ins.InsertAfter($text)
    return $ins
}

# ---------------------------------------------------------------------------
# 1) "geographic boundaries of the UK Local governments [3]" -> "... [11]"
# ---------------------------------------------------------------------------
FindReplace "geographic boundaries of the UK Local governments [3]" "geographic boundaries of the UK Local governments [11]" | Out-Null

Write-Host "Step 1 done"

# ---------------------------------------------------------------------------
# 2) "Properties of the Data" paragraph gets rewritten/expanded into four
#    paragraphs giving more detail about each of the three datasets used.
# ---------------------------------------------------------------------------

# 2a) Rewrite the opening sentence about the COVID case data into the new,
#     much more detailed, first paragraph.
$old2a = "Our COVID 19 case data is at the local authority level (LTLA). The data we have from the census is at the same level, but some of the councils have been merged or split apart."
$new2a = "Our COVID-19 case data [2] is at the local authority level (LTLA) from the results of PCR tests and positive lateral flow tests (which are reported, from 21/10/20). This data is collected from the various local authorities and then checked and published by Public Health England. The data itself is 8 columns by 244,442 rows from 13/3/20 to 28/12/21, with at least one row per day. There is an issue with the data for 1/7/20, which appears to be a correction for earlier data points. Generally, the aggregated local authority data matches the UK wide data, but does not face the same level of scrutiny as the UK wide data (being presented by the Prime Minister)."
FindReplace $old2a $new2a | Out-Null

Write-Host "Step 2a done"

# 2b) Split into a new paragraph right before "Using Excel [4]" and give the
#     new (second) paragraph its opening sentence about the census data.
InsertParagraphAfterAnchor "(being presented by the Prime Minister)." "The data we have from the census [1] is at the same level, but some of the councils have been merged or split apart." | Out-Null

Write-Host "Step 2b done"

# 2c) Append extra sentences about how the census data was collected to the
#     end of the (still same) second paragraph, right after
#     "tables into one sheet."
$r2c = $d.Content
$ok2c = $r2c.Find.Execute("tables into one sheet.")
if (-not $ok2c) { throw "Step 2c anchor not found" }
$r2c.Collapse(0)
$r2c.InsertAfter(" The census data was collected through questionnaires presented to every household in the UK. These were then aggregated by the Office of National Statistics. This dataset has 95 columns and 343 rows.")

Write-Host "Step 2c done"

# 2d) New third paragraph: the local authority boundary data (2020, ONS).
InsertParagraphAfterAnchor "This dataset has 95 columns and 343 rows." "Finally, we have the local authority boundary data for 2020 from [11], which is from the ONS. This dataset has the various local authorities and their geographic properties." | Out-Null

Write-Host "Step 2d done"

# 2e) New fourth paragraph: the three datasets are joined by geography code.
InsertParagraphAfterAnchor "their geographic properties." "The three datasets are joined together by their geography code. By plotting the cases we could see any obvious anomalies in the case data." | Out-Null

Write-Host "Step 2e done"
